# Burndown chart update: record 1 unit of completed effort on "Day 7" (column L)
# for tasks 7 and 8 (rows 12 and 13).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

$ws.Range("L12").Value = 1
$ws.Range("L13").Value = 1

# Leave the selection where the user's last edit was, matching the saved file.
$null = $ws.Range("L13").Select()
